$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Table header: merge "(Sessione " + "Ottobre" + " 2024)" (with
#    proofErr wrapper around "Ottobre") into a single run/text.
# ---------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "(Sessione Ottobre 2024)", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Sessione Ottobre 2024)", 2)
Write-Output "step1 (Sessione Ottobre 2024) merge: $found"

# ---------------------------------------------------------------
# 2) Resize first picture (Diagnosis Flow), width 428.25pt -> 427.9pt.
#    The VML <w:pict> image isn't exposed as an InlineShape/Shape by
#    this runtime, so locate its (otherwise-empty) paragraph via the
#    caption that follows it and rewrite that paragraph's XML with
#    Range.InsertXML (only on the exact range whose content changes).
# ---------------------------------------------------------------
$count = $d.Paragraphs.Count
$capIdx1 = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("Figure 1 - Diagnosis flow")) {
        $capIdx1 = $i
    }
}
Write-Output "step2 caption1 index: $capIdx1"

$picPara1 = $d.Paragraphs.Item($capIdx1 - 1)
$picRange1 = $picPara1.Range
Write-Output "step2 picRange1 text=[$($picRange1.Text)] start=$($picRange1.Start) end=$($picRange1.End)"

$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="6AB4ABBC" w14:textId="5D95824D" w:rsidR="00D33864" w:rsidRDefault="009B058E" w:rsidP="00D33864" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:keepNext/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:pict w14:anchorId="43320D92" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:427.9pt;height:194.25pt"><v:imagedata r:id="rId6" o:title="IMG_1_Diagnosis_Flow"/></v:shape></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$picRange1.InsertXML($xmlFrag1)
Write-Output "step2 done"

# ---------------------------------------------------------------
# 3) Merge the three "ADIA (from Automatic " / "DIAgnostic" / ")..."
#    runs (with proofErr spell-check wrappers) back into a single run.
# ---------------------------------------------------------------
$oldAdia = "The main component of the diagnostic flow is the filtering and validation finite state machine (FSM), called ADIA (from Automatic DIAgnostic). For each fault defined in the system, an automaton cycles through its states to validate (or de-validate) the error conditions identified by the dedicated, user-defined fault check test routine. The ADIA is the core of the diagnostic system, as it" + [char]0x2019 + "s in charge of triggering the faults" + [char]0x2019 + " memory management strategies and preparing the field for the recovery strategies. The original code defined multiple types of ADIAs that have been now condensed into a simpler, standard-compliant single type of FSM ("
$found3 = $d.Content.Find.Execute(
    $oldAdia, $true, $false, $false, $false, $false,
    $true, 1, $false, $oldAdia, 2)
Write-Output "step3 ADIA merge: $found3"

# ---------------------------------------------------------------
# 4) Resize second picture (new ADIA simplified), width 523.5pt ->
#    523.7pt, height 251.25pt -> 251.3pt. Same InsertXML technique.
# ---------------------------------------------------------------
$count = $d.Paragraphs.Count
$capIdx2 = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("New FSM scheme")) {
        $capIdx2 = $i
    }
}
Write-Output "step4 caption2 index: $capIdx2"

$picPara2 = $d.Paragraphs.Item($capIdx2 - 1)
$picRange2 = $picPara2.Range
Write-Output "step4 picRange2 text=[$($picRange2.Text)] start=$($picRange2.Start) end=$($picRange2.End)"

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="63E29448" w14:textId="77777777" w:rsidR="00381BEC" w:rsidRDefault="009B058E" w:rsidP="00381BEC" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:keepNext/><w:jc w:val="center"/></w:pPr><w:r><w:pict w14:anchorId="59D7BAFA" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:523.7pt;height:251.3pt"><v:imagedata r:id="rId7" o:title="IMG_2_new_ADIA_SIMPLIFIED"/></v:shape></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$picRange2.InsertXML($xmlFrag2)
Write-Output "step4 done"

# ---------------------------------------------------------------
# 5) Split " - New FSM scheme, extremely simplified for graphical
#    reasons" into " - New FSM scheme, extremely simplified " (same
#    run) + a new run "to protect the actual know-how" that carries
#    identical formatting (italic, color, Times New Roman incl. cs).
# ---------------------------------------------------------------
$found5a = $d.Content.Find.Execute(
    " - New FSM scheme, extremely simplified for graphical reasons",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " - New FSM scheme, extremely simplified to protect the actual know-how", 2)
Write-Output "step5a combine: $found5a"

$r5 = $d.Content
$found5b = $r5.Find.Execute(
    "to protect the actual know-how", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
Write-Output "step5b find: $found5b text=[$($r5.Text)]"

# Toggling Bold true->false on this exact sub-range forces the engine
# to split it into its own run while it still carries every formatting
# attribute (rFonts incl. cs, bCs, iCs, color) inherited from the
# original run.
$r5.Font.Bold = $true
$r5.Font.Bold = $false
Write-Output "step5 done"
